$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target dataset for columns A, B, C starting at row 2 (row 1 is header).
# Built as an explicit object[] of object[] rows so the jagged structure
# survives (a plain nested @(@(...),@(...)) literal gets flattened).
$data = New-Object 'object[]' 68
$data[0] = @("Alemania (GER)", "GER", "Alemania")
$data[1] = @("Argelia", $null, "Argelia")
$data[2] = @("Argentina (ARG)", "ARG", "Argentina")
$data[3] = @("Armenia (ARM)", "ARM", "Armenia")
$data[4] = @("Australia (AUS)", "AUS", "Australia")
$data[5] = @("Austria (AUT)", "AUT", "Austria")
$data[6] = @("Azerbaiyán (AZE)", "AZE", "Azerbaiyán")
$data[7] = @("Brasil (BRA)", "BRA", "Brasil")
$data[8] = @("Bélgica (BEL)", "BEL", "Bélgica")
$data[9] = @("Cabo Verde", $null, "Cabo Verde")
$data[10] = @("Canadá (CAN)", "CAN", "Canadá")
$data[11] = @("Chile", $null, "Chile")
$data[12] = @("China (CHN)", "CHN", "China")
$data[13] = @("Corea del Norte (PRK)", "PRK", "Corea del Norte")
$data[14] = @("Corea del Sur (KOR)", "KOR", "Corea del Sur")
$data[15] = @("Croacia (CRO)", "CRO", "Croacia")
$data[16] = @("Cuba (CUB)", "CUB", "Cuba")
$data[17] = @("Dinamarca (DEN)", "DEN", "Dinamarca")
$data[18] = @("Dominica", $null, "Dominica")
$data[19] = @("Ecuador (ECU)", "ECU", "Ecuador")
$data[20] = @("Egipto (EGY)", "EGY", "Egipto")
$data[21] = @("Eslovaquia (SVK)", "SVK", "Eslovaquia")
$data[22] = @("Eslovenia (SLO)", "SLO", "Eslovenia")
$data[23] = @("España (ESP)", "ESP", "España")
$data[24] = @("Estados Unidos (USA)", "USA", "Estados Unidos")
$data[25] = @("Etiopía (ETH)", "ETH", "Etiopía")
$data[26] = @("Filipinas (PHI)", "PHI", "Filipinas")
$data[27] = @("Fiyi (FIY)", "FIY", "Fiyi")
$data[28] = @("Francia (FRA)", "FRA", "Francia")
$data[29] = @("Georgia (GEO)", "GEO", "Georgia")
$data[30] = @("Grecia (GRE)", "GRE", "Grecia")
$data[31] = @("Guatemala (GUA)", "GUA", "Guatemala")
$data[32] = @("Hong Kong (HKG)", "HKG", "Hong Kong")
$data[33] = @("Hungría (HUN)", "HUN", "Hungría")
$data[34] = @("India (IND)", "IND", "India")
$data[35] = @("Indonesia (INA)", "INA", "Indonesia")
$data[36] = @("Irlanda (IRL)", "IRL", "Irlanda")
$data[37] = @("Israel (ISR)", "ISR", "Israel")
$data[38] = @("Italia (ITA)", "ITA", "Italia")
$data[39] = @("Jamaica (JAM)", "JAM", "Jamaica")
$data[40] = @("Japón (JPN)", "JPN", "Japón")
$data[41] = @("Kazajistán (KAZ)", "KAZ", "Kazajistán")
$data[42] = @("Kosovo (KOS)", "KOS", "Kosovo")
$data[43] = @("Lituania (LTU)", "LTU", "Lituania")
$data[44] = @("Malasia (MAL)", "MAL", "Malasia")
$data[45] = @("Moldavia (MDA)", "MDA", "Moldavia")
$data[46] = @("Mongolia (MGL)", "MGL", "Mongolia")
$data[47] = @("México (MEX)", "MEX", "México")
$data[48] = @("Nueva Zelanda (NZL)", "NZL", "Nueva Zelanda")
$data[49] = @("Países Bajos (NED)", "NED", "Países Bajos")
$data[50] = @("Polonia (POL)", "POL", "Polonia")
$data[51] = @("Portugal (POR)", "POR", "Portugal")
$data[52] = @("Reino Unido (GBR)", "GBR", "Reino Unido")
$data[53] = @("República Checa (CZE)", "CZE", "República Checa")
$data[54] = @("República Dominicana (DOM)", "DOM", "República Dominicana")
$data[55] = @("Rumania (ROU)", "ROU", "Rumania")
$data[56] = @("Santa Lucía", $null, "Santa Lucía")
$data[57] = @("Serbia (SRB)", "SRB", "Serbia")
$data[58] = @("Sudáfrica (RSA)", "RSA", "Sudáfrica")
$data[59] = @("Suecia (SWE)", "SWE", "Suecia")
$data[60] = @("Suiza (SUI)", "SUI", "Suiza")
$data[61] = @("Taiwán", $null, "Taiwán")
$data[62] = @("Tayikistán (TJK)", "TJK", "Tayikistán")
$data[63] = @("Turquía (TUR)", "TUR", "Turquía")
$data[64] = @("Túnez (TUN)", "TUN", "Túnez")
$data[65] = @("Ucrania (UKR)", "UKR", "Ucrania")
$data[66] = @("Uganda (UGA)", "UGA", "Uganda")
$data[67] = @("Uzbekistán (UZB)", "UZB", "Uzbekistán")

$oldLastRow = 51

# Clear out the previous data body (rows 2..51) in columns A:C so stale cells
# do not linger if the new table ends up shorter than before.
$clearRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($oldLastRow, 3))
$clearRange.ClearContents()

# Write the new rows one by one so that countries without an ISO code leave
# column B blank instead of inheriting a stale value.
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]

    $ws.Cells.Item($rowNum, 1).Value = $row[0]
    if ($null -eq $row[1]) {
        $ws.Cells.Item($rowNum, 2).Value = ""
    } else {
        $ws.Cells.Item($rowNum, 2).Value = $row[1]
    }
    $ws.Cells.Item($rowNum, 3).Value = $row[2]
}
